$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced
# to text so Excel keeps them as strings (matching the source data format)
# instead of auto-converting them to numeric cells.
$ws.Range("D2").Value = "69.477.44"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.503.00"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -6.37%  "
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "4.054.78"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "597.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.81%  "
$ws.Range("D16").Value = "69.615.45"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.123"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "3.500.18"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.77%  "
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("E30").Value = "  +9.01%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").Value = "3.746.12"
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("D37").Value = "0.0₃0811"
$ws.Range("E37").Value = "  +4.82%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "492.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.93%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0449"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.39%  "
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("B50").Value = "OceanProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000243"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.97%  "
